$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# ---- Add the three new sheets in order, right after Sheet1 ----
$wsUser = $wb.Worksheets.Add($null, $sheet1)
$wsUser.Name = "User"

$wsCustomers = $wb.Worksheets.Add($null, $wsUser)
$wsCustomers.Name = "Customers"

$wsCustomersLocation = $wb.Worksheets.Add($null, $wsCustomers)
$wsCustomersLocation.Name = "CustomersLocation"

# ---- Sheet1: view scrolled down a bit, no longer the tab shown on open ----
$sheet1.Range("A1").Select()
$sheet1.Application.ActiveWindow.ScrollRow = 7

# ================= User sheet =================
$wsUser.Range("A2").Value = "ID"
$wsUser.Range("B2").Value = "int"
$wsUser.Range("D2").Value = "Mô tả"

$wsUser.Range("A3").Value = "UserName"
$wsUser.Range("B3").Value = "nvarchar(50)"

$wsUser.Range("A4").Value = "Password"
$wsUser.Range("B4").Value = "nvarchar(50)"

$wsUser.Range("A5").Value = "Op"
$wsUser.Range("B5").Value = "int"
$wsUser.Range("D5").Value = "Nếu OP = 0 thì Sau khi đăng nhập sẽ vào giao diện chính, Op=1 Thì sau khi đăng nhập sẽ ra bảng lọc khách hàng. Sau khi lọc khách hàng xong thì chọn vào 1 khách hàng  sẽ lưu lại tọa độ và vào giao diện chính(tất cả các chức năng phía sau sẽ ăn theo khách hàng vừa chọn)"

$wsUser.Range("A6").Value = "IdenUser"
$wsUser.Range("B6").Value = "int"
$wsUser.Range("D6").Value = "Chứa Identity ngầm của phần mềm"

$wsUser.Range("A7").Value = "DateTime"
$wsUser.Range("B7").Value = "datetime"

# Styling: rows 3-6 col A:B -> black font, left/vcenter aligned
$wsUser.Range("A3:B6").Font.Color = 0
$wsUser.Range("A3:B6").HorizontalAlignment = -4131
$wsUser.Range("A3:B6").VerticalAlignment = -4108
# Row 7 col A:B -> black font only (no alignment override)
$wsUser.Range("A7:B7").Font.Color = 0
# D5 description cell -> default font, left aligned, taller row
$wsUser.Range("D5").HorizontalAlignment = -4131
$wsUser.Rows("5").RowHeight = 30.75

$wsUser.Columns("A").ColumnWidth = 16.140625
$wsUser.Columns("B").ColumnWidth = 13.28515625
$wsUser.Columns("D").ColumnWidth = 86.140625

$wsUser.Range("E7").Select()

# ================= Customers sheet =================
$wsCustomers.Range("A2").Value = "Các trường"
$wsCustomers.Range("B2").Value = "Kiểu dữ liệu"
$wsCustomers.Range("D2").Value = "Mô tả"

$wsCustomers.Range("A3").Value = "ID"
$wsCustomers.Range("B3").Value = "int"

$wsCustomers.Range("A4").Value = "CustomersId"
$wsCustomers.Range("B4").Value = "nvarchar(50)"
$wsCustomers.Range("D4").Value = "Chứa mã khách hàng"

$wsCustomers.Range("A5").Value = "CustomersName"
$wsCustomers.Range("B5").Value = "nvarchar(100)"
$wsCustomers.Range("D5").Value = "Chứa tên khách hàng"

$wsCustomers.Range("A6").Value = "Adress"
$wsCustomers.Range("B6").Value = "nvarhar(200)"
$wsCustomers.Range("D6").Value = "Địa chỉ khách hàng"

$wsCustomers.Range("A7").Value = "Phone"
$wsCustomers.Range("B7").Value = "nvarchar(30)"
$wsCustomers.Range("D7").Value = "Số điện thoại"

$wsCustomers.Range("A8").Value = "IdenCustomers"
$wsCustomers.Range("B8").Value = "int"
$wsCustomers.Range("D8").Value = "Chứa Identity của phần mềm tương ứng với bảng khách hàng"

$wsCustomers.Range("A9").Value = "IdenUser"
$wsCustomers.Range("B9").Value = "int"
$wsCustomers.Range("D9").Value = "Chứa Identity của phần mềm tương ứng với bảng User"

$wsCustomers.Range("A10").Value = "DateTime"
$wsCustomers.Range("B10").Value = "datetime"

$wsCustomers.Range("A3:B9").Font.Color = 0
$wsCustomers.Range("A3:B9").HorizontalAlignment = -4131
$wsCustomers.Range("A3:B9").VerticalAlignment = -4108
$wsCustomers.Range("A10:B10").Font.Color = 0

$wsCustomers.Columns("A").ColumnWidth = 19.42578125
$wsCustomers.Columns("B").ColumnWidth = 17.85546875

$wsCustomers.Range("G11").Select()

# ================= CustomersLocation sheet =================
$wsCustomersLocation.Range("A2").Value = "Các trường"
$wsCustomersLocation.Range("B2").Value = "Kiểu dữ liệu"
$wsCustomersLocation.Range("D2").Value = "Mô tả"

$wsCustomersLocation.Range("A3").Value = "ID"
$wsCustomersLocation.Range("B3").Value = "int"

$wsCustomersLocation.Range("A4").Value = "IdenCustomers"
$wsCustomersLocation.Range("B4").Value = "int"
$wsCustomersLocation.Range("D4").Value = "Chứa Identity của bảng khách hàng"

$wsCustomersLocation.Range("A5").Value = "LocationName"
$wsCustomersLocation.Range("B5").Value = "nvarchar(500)"
$wsCustomersLocation.Range("D5").Value = "Chứa tên địa điểm"

$wsCustomersLocation.Range("A6").Value = "Description"
$wsCustomersLocation.Range("B6").Value = "nvarchar(500)"
$wsCustomersLocation.Range("D6").Value = "Chứa Nội dung địa điểm"

$wsCustomersLocation.Range("A7").Value = "longitude"
$wsCustomersLocation.Range("B7").Value = "nvarchar(200)"
$wsCustomersLocation.Range("D7").Value = "Chứa kinh độ "

$wsCustomersLocation.Range("A8").Value = "Latitude"
$wsCustomersLocation.Range("B8").Value = "nvarchar(200)"
$wsCustomersLocation.Range("D8").Value = "Chứa vĩ độ"

$wsCustomersLocation.Range("A9").Value = "IdenUser"
$wsCustomersLocation.Range("B9").Value = "int"
$wsCustomersLocation.Range("D9").Value = "Chứa identity của bảng user (iden phần mềm)"

$wsCustomersLocation.Range("A10").Value = "IdenCustomersLocation"
$wsCustomersLocation.Range("B10").Value = "int"
$wsCustomersLocation.Range("D10").Value = "Chứa iden để đồng bộ"

$wsCustomersLocation.Range("A11").Value = "DateTime"
$wsCustomersLocation.Range("B11").Value = "datetime"

$wsCustomersLocation.Range("A12").Value = "Disable"
$wsCustomersLocation.Range("B12").Value = "int"
$wsCustomersLocation.Range("D12").Value = "nếu =1 ko hiển thị"

$wsCustomersLocation.Range("H16").Value = "Haithanh#1"
$wsCustomersLocation.Range("H18").Value = "TTDH"
$wsCustomersLocation.Range("H19").Value = "1+2=3"

$wsCustomersLocation.Range("A3:B11").Font.Color = 0
$wsCustomersLocation.Range("A3:B11").HorizontalAlignment = -4131
$wsCustomersLocation.Range("A3:B11").VerticalAlignment = -4108
$wsCustomersLocation.Range("A12:B12").Font.Color = 0

$wsCustomersLocation.Range("A2:E2").Font.Bold = $true
$wsCustomersLocation.Range("A2:E2").Font.Color = 255
$wsCustomersLocation.Range("A2:E2").Interior.Color = 65535

$wsCustomersLocation.Columns("A").ColumnWidth = 18.42578125
$wsCustomersLocation.Columns("B").ColumnWidth = 14.5703125
$wsCustomersLocation.Columns("G").ColumnWidth = 33.42578125
$wsCustomersLocation.Columns("H").ColumnWidth = 39.85546875

$wsCustomersLocation.Range("H18").Select()

# ---- Make the User sheet the active/shown tab, matching activeTab=1 ----
$wsUser.Activate()
